# Sync "Recorded By" (column G) attendance-report values: for any cell whose
# comma-separated list starts with "System"/"system", re-sort the whole list
# alphabetically (case-insensitive) - matching upstream main repo sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value()
    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq "system") {
        $sorted = $parts | Sort-Object { $_.ToLower() }
        $newText = [string]::Join(", ", $sorted)
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
